$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.78
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.35
$ws.Range("K3").Value = 2.02
$ws.Range("T3").Value = 2.55
$ws.Range("U3").Value = 1.98
$ws.Range("V3").Value = 1.75
$ws.Range("W3").Value = 5.8
$ws.Range("X3").Value = 7.7
$ws.Range("Z3").Value = 14.5
$ws.Range("AH3").Value = 11
$ws.Range("AI3").Value = 28
$ws.Range("AM3").Value = 65
$ws.Range("AN3").Value = 3.55
$ws.Range("AT3").Value = 2.55
$ws.Range("AV3").Value = 75
$ws.Range("AW3").Value = 6.6
$ws.Range("AX3").Value = 32
$ws.Range("AY3").Value = 37
$ws.Range("BB3").Value = 500

# Row 5 updates
$ws.Range("G5").Value = 1.33
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 1.8
$ws.Range("K5").Value = 2.4
$ws.Range("L5").Value = 7.5
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.08
$ws.Range("S5").Value = 1.33
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("X5").Value = 6
$ws.Range("Z5").Value = 8.5
$ws.Range("AA5").Value = 12
$ws.Range("AC5").Value = 13
$ws.Range("AD5").Value = 10
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 17
$ws.Range("AI5").Value = 41
$ws.Range("AJ5").Value = 23
$ws.Range("AN5").Value = 3.25
$ws.Range("AO5").Value = 6.5
$ws.Range("AQ5").Value = 17
$ws.Range("AT5").Value = 3.25
$ws.Range("AU5").Value = 10
$ws.Range("AW5").Value = 9
